$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 709; this shifts the existing rows 709-793 down to 710-794
$ws.Rows("709:709").Insert()

# Populate the newly inserted row 709 with the new weekly record
$ws.Range("A709").Value = 7
$ws.Range("B709").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C709").Value = "Ñuble"
$ws.Range("D709").Value = 45212
$ws.Range("E709").Value = 16
$ws.Range("F709").Value = 100114014
$ws.Range("G709").Value = "Betarraga"
$ws.Range("H709").Value = "Sin especificar"
$ws.Range("I709").Value = "Segunda"
$ws.Range("J709").Value = 300
$ws.Range("K709").Value = 600
$ws.Range("L709").Value = 600
$ws.Range("M709").Value = 600
$ws.Range("N709").Value = "$/paquete 5 unidades"
$ws.Range("O709").Value = "Provincia de Diguillín"
$ws.Range("P709").Value = 120
$ws.Range("Q709").Value = 5
$ws.Range("R709").Value = "Hortaliza"
